$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Наличные")

$ws.Cells.Item(35, 1).Value = 7676096317
$ws.Cells.Item(35, 2).Value = "M-Банкинг чек-4294968789.docx"
$ws.Cells.Item(35, 3).Value = 1
$ws.Cells.Item(35, 4).Value = 0.2
$ws.Cells.Item(35, 5).Value = "2025-06-28 01:33:45"
